$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 271, shifting existing rows 271:330 down to 272:331
$ws.Rows("271:271").Insert()

$ws.Range("A272").Copy()
$ws.Range("A271").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A271").Value = "FRU"
$ws.Range("B271").Value = "Bishkek, Kyrgyzstan"
$ws.Range("C271").Value = "Asia Pacific"
$ws.Range("D271").Value = "Bishkek"
$ws.Range("E271").Value = "Kyrgyzstan"
$ws.Range("F271").Value = "KG"
$ws.Range("G271").Value = 42.875608
$ws.Range("H271").Value = 74.604613
